$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.487.89'
$ws.Range("E2").Value = '  +0.78%  '
$ws.Range("E3").Value = '  +1.18%  '
$ws.Range("E4").Value = '  +0.21%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '219.63'
$ws.Range("E5").Value = '  +0.68%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5496'
$ws.Range("E6").Value = '  +4.85%  '
$ws.Range("E7").Value = '  +0.16%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2746'
$ws.Range("E8").Value = '  +1.38%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06484'
$ws.Range("E9").Value = '  +1.05%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.13'
$ws.Range("E10").Value = '  +0.65%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07706'
$ws.Range("E11").Value = '  +2.72%  '
$ws.Range("D12").Value = '1.728.48'
$ws.Range("E12").Value = '  +2.17%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.565'
$ws.Range("E13").Value = '  +0.10%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5855'
$ws.Range("E14").Value = '  +1.09%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.000008428'
$ws.Range("E15").Value = '  -0.26%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.22'
$ws.Range("E16").Value = '  +3.07%  '
$ws.Range("D17").Value = '26.563.40'
$ws.Range("E17").Value = '  +0.89%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.970'
$ws.Range("E18").Value = '  +1.03%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.010'
$ws.Range("E19").Value = '  +0.15%  '
$ws.Range("E20").Value = '  +1.37%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '192.53'
$ws.Range("E21").Value = '  +2.13%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.284'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.011'
$ws.Range("E23").Value = '  +0.15%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '149.47'
$ws.Range("E24").Value = '  +3.46%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1331'
$ws.Range("E25").Value = '  +8.00%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.944'
$ws.Range("E26").Value = '  +3.20%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.87'
$ws.Range("E27").Value = '  +0.45%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.06294'
$ws.Range("E28").Value = '  -5.65%  '
$ws.Range("E29").Value = '  +2.57%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.334'
$ws.Range("E30").Value = '  +0.46%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.630'
$ws.Range("E31").Value = '  +1.77%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.613'
$ws.Range("E32").Value = '  +1.03%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.698'
$ws.Range("E33").Value = '  +2.39%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.049'
$ws.Range("E34").Value = '  +2.07%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6209'
$ws.Range("E35").Value = '  +0.15%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.412'
$ws.Range("E36").Value = '  +0.53%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.774'
$ws.Range("E37").Value = '  +2.99%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01651'
$ws.Range("E38").Value = '  +2.19%  '
$ws.Range("D39").Value = '1.124.01'
$ws.Range("E39").Value = '  +1.71%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.164'
$ws.Range("E40").Value = '  -3.48%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8852'
$ws.Range("E41").Value = '  +0.98%  '
$ws.Range("E42").Value = '  +0.26%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '101.59'
$ws.Range("E43").Value = '  +0.85%  '
$ws.Range("D44").Value = '1.855.16'
$ws.Range("E44").Value = '  +1.24%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '57.91'
$ws.Range("E45").Value = '  +2.13%  '
$ws.Range("E46").Value = '  -1.06%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.291'
$ws.Range("E47").Value = '  +1.22%  '
$ws.Range("E48").Value = '  -0.43%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05286'
$ws.Range("E49").Value = '  +0.26%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.157'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4305'
$ws.Range("E51").Value = '  -0.02%  '
